$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 527.75
$ws.Range("I19").Value = 360.25
$ws.Range("J19").Value = 594.75
$ws.Range("K19").Value = 360.25
$ws.Range("L19").Value = 594.75
$ws.Range("M19").Value = -185.25
$ws.Range("N19").Value = -944.75
$ws.Range("H20").Value = 595
$ws.Range("I20").Value = 595
$ws.Range("K20").Value = 595
$ws.Range("M20").Value = -365
$ws.Range("H35").Value = 595
$ws.Range("I35").Value = 595
$ws.Range("K35").Value = 595
$ws.Range("M35").Value = -216
$ws.Range("H98").Value = 1524.1305
$ws.Range("I98").Value = 1132.3529
$ws.Range("J98").Value = 2634.1667
$ws.Range("K98").Value = 1132.3529
$ws.Range("L98").Value = 2634.1667
$ws.Range("M98").Value = 365.6470999999999
$ws.Range("N98").Value = -5630.1667
$ws.Range("H122").Value = 1524.1305
$ws.Range("I122").Value = 1132.3529
$ws.Range("J122").Value = 2634.1667
$ws.Range("K122").Value = 3397.0587
$ws.Range("L122").Value = 7902.500100000001
$ws.Range("M122").Value = -947.0587000000005
$ws.Range("N122").Value = -12802.5001
$ws.Range("H129").Value = 1347.9286
$ws.Range("I129").Value = 479.2857
$ws.Range("J129").Value = 1637.4762
$ws.Range("K129").Value = 1437.8571
$ws.Range("L129").Value = 4912.4286
$ws.Range("M129").Value = 3562.1429
$ws.Range("N129").Value = -14912.4286
$ws.Range("H132").Value = 6849.9316
$ws.Range("I132").Value = 5458.7354
$ws.Range("J132").Value = 11580
$ws.Range("K132").Value = 16376.2062
$ws.Range("L132").Value = 34740
$ws.Range("M132").Value = -13846.2062
$ws.Range("N132").Value = -39800
$ws.Range("H137").Value = 8476362
$ws.Range("I137").Value = 13159728
$ws.Range("J137").Value = 1698.8096
$ws.Range("K137").Value = 39479184
$ws.Range("L137").Value = 5096.4288
$ws.Range("M137").Value = -39476634
$ws.Range("N137").Value = -10196.4288

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4515609.5
$ws.Range("I32").Value = 6043.212
$ws.Range("K32").Value = 6043.212
$ws.Range("M32").Value = -5756.212
$ws.Range("H61").Value = 2896.7932
$ws.Range("I61").Value = 1727.6666
$ws.Range("J61").Value = 3722.0588
$ws.Range("K61").Value = 1727.6666
$ws.Range("L61").Value = 3722.0588
$ws.Range("M61").Value = -1515.6666
$ws.Range("N61").Value = -4146.0588
$ws.Range("H136").Value = 2896.7932
$ws.Range("I136").Value = 1727.6666
$ws.Range("J136").Value = 3722.0588
$ws.Range("K136").Value = 5182.9998
$ws.Range("L136").Value = 11166.1764
$ws.Range("M136").Value = -2632.9998
$ws.Range("N136").Value = -16266.1764

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 640.5714
$ws.Range("I94").Value = 580.6667
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 580.6667
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -129.6667
$ws.Range("N94").Value = -1902

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1434738.1
$ws.Range("I58").Value = 1393.65
$ws.Range("J58").Value = 3345864.2
$ws.Range("K58").Value = 1393.65
$ws.Range("L58").Value = 3345864.2
$ws.Range("M58").Value = -1190.65
$ws.Range("N58").Value = -3346270.2
$ws.Range("H64").Value = 19983.166
$ws.Range("J64").Value = 19983.166
$ws.Range("L64").Value = 19983.166
$ws.Range("N64").Value = -20479.166
$ws.Range("H67").Value = 19983.166
$ws.Range("J67").Value = 19983.166
$ws.Range("L67").Value = 19983.166
$ws.Range("N67").Value = -21699.166
$ws.Range("H132").Value = 5540.8887
$ws.Range("I132").Value = 4467.5
$ws.Range("J132").Value = 6399.6
$ws.Range("K132").Value = 13402.5
$ws.Range("L132").Value = 19198.8
$ws.Range("M132").Value = -10872.5
$ws.Range("N132").Value = -24258.8
$ws.Range("H134").Value = 3642.7856
$ws.Range("I134").Value = 1872.75
$ws.Range("J134").Value = 4350.8
$ws.Range("K134").Value = 5618.25
$ws.Range("L134").Value = 13052.4
$ws.Range("M134").Value = -3083.25
$ws.Range("N134").Value = -18122.4
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 1434738.1
$ws.Range("I136").Value = 1393.65
$ws.Range("J136").Value = 3345864.2
$ws.Range("K136").Value = 4180.950000000001
$ws.Range("L136").Value = 10037592.6
$ws.Range("M136").Value = -1630.950000000001
$ws.Range("N136").Value = -10042692.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 875272.5
$ws.Range("I68").Value = 1027.0769
$ws.Range("J68").Value = 1749517.9
$ws.Range("K68").Value = 3081.2307
$ws.Range("L68").Value = 5248553.699999999
$ws.Range("M68").Value = -2270.2307
$ws.Range("N68").Value = -5250175.699999999
$ws.Range("H71").Value = 875272.5
$ws.Range("I71").Value = 1027.0769
$ws.Range("J71").Value = 1749517.9
$ws.Range("K71").Value = 9243.6921
$ws.Range("L71").Value = 15745661.1
$ws.Range("M71").Value = -5187.6921
$ws.Range("N71").Value = -15753773.1
$ws.Range("H107").Value = 15955163
$ws.Range("J107").Value = 619979.6
$ws.Range("L107").Value = 1859938.8
$ws.Range("N107").Value = -1863778.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2491.7273
$ws.Range("I7").Value = 2763
$ws.Range("J7").Value = 1768.3334
$ws.Range("K7").Value = 2763
$ws.Range("L7").Value = 1768.3334
$ws.Range("M7").Value = -2651
$ws.Range("N7").Value = -1992.3334
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H122").Value = 9440
$ws.Range("I122").Value = 13675.2
$ws.Range("J122").Value = 5910.6665
$ws.Range("K122").Value = 41025.60000000001
$ws.Range("L122").Value = 17731.9995
$ws.Range("M122").Value = -38575.60000000001
$ws.Range("N122").Value = -22631.9995
$ws.Range("H126").Value = 2491.7273
$ws.Range("I126").Value = 2763
$ws.Range("J126").Value = 1768.3334
$ws.Range("K126").Value = 8289
$ws.Range("L126").Value = 5305.0002
$ws.Range("M126").Value = -5819
$ws.Range("N126").Value = -10245.0002

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 7584653
$ws.Range("I136").Value = 12513633
$ws.Range("K136").Value = 37540899
$ws.Range("M136").Value = -37538349

Write-Output "Applied all changes"